$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("baiama",    10, 1, 1,   49),
    @("baiama",    10, 2, 50,  49),
    @("baiama",    10, 3, 99,  49),
    @("baiama",    10, 4, 148, 49),
    @("baiama",    10, 7, 197, 98),
    @("lalehun",   10, 1, 1,   49),
    @("lalehun",   10, 2, 50,  49),
    @("lalehun",   10, 3, 99,  49),
    @("lalehun",   10, 4, 147, 49),
    @("lalehun",   10, 5, 197, 49),
    @("lalehun",   10, 7, 246, 94),
    @("lambayama", 10, 1, 1,   49),
    @("lambayama", 10, 2, 50,  49),
    @("lambayama", 10, 3, 99,  49),
    @("lambayama", 10, 4, 148, 49),
    @("lambayama", 10, 7, 197, 98),
    @("seilama",   10, 1, 1,   49),
    @("seilama",   10, 2, 50,  49),
    @("seilama",   10, 3, 99,  49),
    @("seilama",   10, 4, 148, 49),
    @("seilama",   10, 5, 197, 49),
    @("seilama",   10, 7, 246, 98)
)

$startRow = 17
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

$ws.Range("E38").Select()
